# Rename the worksheet from "AlphaFiberF-HW45.xpc" to "AlphaFiberF"
# and append a new data row (row 16) for the "HexGrid-60degTilt5degRes"
# Gaussian-quadrature scheme, matching the existing rows' pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename sheet
$ws.Name = "AlphaFiberF"

# 2) Append new row 16: index 14 in column A, scheme name in column B,
#    and 1s across C:M (same shape as every other data row).
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16:M16").Value = 1

# 3) Match formatting of column A's existing header cells (bold, bordered,
#    centered style) by copying the format from the row above.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
